# Branch wise stock status added
# Re-order a handful of item rows (Dinafex / Ketonic / Kynol / Zithrox groups)
# so that the Item Name (D) / UOM (E) pairing for each row reflects the new
# intended ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dinafex group: swap 120mg and 60mg rows
$ws.Range("D4").Value = "Dinafex 60mg Tablet"
$ws.Range("D5").Value = "Dinafex 120mg Tablet"

# Ketonic / Kynol group (rows 14-19)
$ws.Range("D14").Value = "Ketonic 10mg Tablet"
$ws.Range("E14").Value = "20's"

$ws.Range("D16").Value = "Ketonic 30mg Injection"
$ws.Range("E16").Value = "5 's"

$ws.Range("D17").Value = "Kynol TR 100mg Capsule"
$ws.Range("E17").Value = "50 's"

$ws.Range("D18").Value = "Kynol TR 200mg Capsule"
$ws.Range("E18").Value = "30 's"

$ws.Range("D19").Value = "Kynol D 25mg Tablet"
$ws.Range("E19").Value = "60 's"

# Zithrox group (rows 24-27)
$ws.Range("D24").Value = "Zithrox 15ml Suspension"
$ws.Range("E24").Value = "15 ml"

$ws.Range("D25").Value = "Zithrox 250mg Tablet - 6's"
$ws.Range("E25").Value = "6's"

$ws.Range("D26").Value = "Zithrox 500mg Tablet"
$ws.Range("E26").Value = "6 's"

$ws.Range("D27").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E27").Value = "30ml"
